$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13 & 14: Owner changed from "RJ" to "AJ" ---
$ws.Range("M13").Value = "AJ"
$ws.Range("M14").Value = "AJ"

# --- Row 15: new requirement HLR_009 ---
$ws.Range("A15").Value = "HLR_009"
$ws.Range("B15").Value = "Order Processing"
$ws.Range("C15").Value = "Show Status"
$ws.Range("D15").Value = "System must add Order Status"
$ws.Range("E15").Value = "Low"
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 'System must add an Order Status "Buy" or "Sell  to the order created. '
$ws.Range("H15").Value = "SOW"
$ws.Range("I15").Value = "Processing"
$ws.Range("J15").Value = "High"
$ws.Range("K15").Value = "HLR_001"
$ws.Range("L15").Value = "In-Progress"
$ws.Range("M15").Value = "AJ"
$ws.Range("N15").Value = "N/A"
$ws.Range("O15").Value = "N"
$ws.Range("P15").Value = "None"

# --- Row 16: new requirement HLR_010 ---
$ws.Range("A16").Value = "HLR_010"
$ws.Range("B16").Value = "Order Creation"
$ws.Range("C16").Value = "Party ID"
$ws.Range("D16").Value = "An admin must be able to add the Party ID"
$ws.Range("E16").Value = "High"
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = "The admin must be able to assign a Party ID to the order"
$ws.Range("H16").Value = "SOW"
$ws.Range("I16").Value = "Registration"
$ws.Range("J16").Value = "High"
$ws.Range("K16").Value = "HLR_001,002"
$ws.Range("L16").Value = "In-Progress"
$ws.Range("M16").Value = "AJ"
$ws.Range("N16").Value = "N/A"
$ws.Range("O16").Value = "N"
$ws.Range("P16").Value = "None"

# --- Row 17: new requirement id only ---
$ws.Range("A17").Value = "HLR_011"

# --- Sheet view: scroll to A1 and select B17 ---
$ws.Activate()
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("B17").Select()
